# Add a new column W ("T21: 9/4/2020") to the COVID19 history sheet,
# mirroring the existing column V (T20) layout/formatting, with updated
# data values (row 7 differs: 265 instead of 226) and a new SUM formula
# in row 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell W1 ---------------------------------------------------
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("W1").Value = "T21: 9/4/2020"

# --- Data cells W2:W19 (copy formatting from column V, then set value) -
$wValues = @{
    2  = 16
    3  = 2
    4  = 20
    5  = 2
    6  = 1
    7  = 265
    8  = 0
    9  = 54
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 4
    15 = 0
    16 = 0
    17 = 8
    18 = 0
    19 = 9
}

foreach ($r in 2..19) {
    $ws.Range("V" + $r).Copy()
    $ws.Range("W" + $r).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("W" + $r).Value = $wValues[$r]
}

# --- Sum row (row 20) ---------------------------------------------------
$ws.Range("V20").Copy()
$ws.Range("W20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("W20").Formula = "=SUM(W2:W19)"

$excel.CutCopyMode = 0

# --- Column width for the new column W (approx. 13.75 chars) ----------
$ws.Columns.Item(23).ColumnWidth = 12.92

# --- Selection/active cell, matching the authored workbook state -------
$ws.Range("W20").Select()
